# Generate Report for Handback
#
# The localization pipeline discovered a handback (.xlf) file for
# 56fb2371-e384-4733-9675-3abb445e8db3.md whose handed-back content is
# *not* the latest version of the source document. Record that on row 7
# of both the "zh-cn" and "de-de" status sheets:
#   - Latest Target File (I7)          -> hyperlink to the md doc
#   - Latest Handback File (J7)        -> the handed-back xlf file name
#   - Latest Handback DateTime (K7)    -> when the handback was produced
#   - Error Detail (P7)                -> explanation: stale handback version
#
# Also widen the "Error Detail" column (P) so the long message is readable.

$wb = $excel.ActiveWorkbook

# ColumnWidth is expressed in characters; Excel pads the stored XML
# <col width> by 5/6 of a character versus the ColumnWidth value, so to
# land on a stored width of exactly 40 we dial in 40 - 5/6.
$targetColumnWidth = 40 - (5 / 6)

$sheetsInfo = @(
    @{
        Name = "zh-cn"
        HandbackFile = "56fb2371-e384-4733-9675-3abb445e8db3.54c9de0d4372e3f6ed50b6887219185175f1edee.zh-cn.xlf"
        HandbackDate = "2016-08-14 00:58:25"
    },
    @{
        Name = "de-de"
        HandbackFile = "56fb2371-e384-4733-9675-3abb445e8db3.54c9de0d4372e3f6ed50b6887219185175f1edee.de-de.xlf"
        HandbackDate = "2016-08-14 00:58:35"
    }
)

$targetMdName = "56fb2371-e384-4733-9675-3abb445e8db3.md"
$latestUrl = "https://github.com/OpenLocalizationTestOrg/oltest/blob/094f880e79721e844dd2af18496de45bbe8f84d8/e2e/56fb2371-e384-4733-9675-3abb445e8db3.md"
$currentUrl = "https://github.com/OpenLocalizationTestOrg/oltest/blob/76321ae9c801e197ca0f17f1b9769d7cbe156d0b/e2e/56fb2371-e384-4733-9675-3abb445e8db3.md"
$errorMessage = "The version of handback file is not the latest, current: " + $currentUrl + ", latest: " + $latestUrl + "."

foreach ($info in $sheetsInfo) {
    $ws = $wb.Worksheets.Item($info.Name)

    # Latest Handback File (J7) + Latest Handback DateTime (K7)
    $ws.Range("J7").Value = $info.HandbackFile
    $ws.Range("K7").Value = $info.HandbackDate

    # Error Detail (P7)
    $ws.Range("P7").Value = $errorMessage

    # Latest Target File (I7) becomes a hyperlink to the (non-latest)
    # target markdown doc; Hyperlinks.Add both sets the cell's display
    # text and applies hyperlink styling.
    $ws.Hyperlinks.Add($ws.Range("I7"), $latestUrl, "", "", $targetMdName)

    # Widen the Error Detail column (P / column 16) to fit the message.
    $ws.Columns.Item(16).ColumnWidth = $targetColumnWidth
}
